$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-14 Monday" "2025-04-15 Tuesday"

Replace-Text "44×37=" "16×72="
Replace-Text "91×53=" "15×62="
Replace-Text "72×86=" "80×90="
Replace-Text "54×91=" "41×22="
Replace-Text "31×30=" "55×65="

Replace-Text "53×78=" "25×88="
Replace-Text "73×89=" "14×36="
Replace-Text "22×34=" "91×96="
Replace-Text "31×64=" "50×93="
Replace-Text "48×87=" "71×68="

Replace-Text "78×41=" "52×44="
Replace-Text "31×14=" "36×12="
Replace-Text "64×73=" "81×38="
Replace-Text "48×95=" "19×69="
Replace-Text "86×14=" "30×39="

Replace-Text "20×85=" "60×40="
Replace-Text "56×95=" "13×35="
Replace-Text "65×40=" "99×89="
Replace-Text "17×62=" "93×46="
Replace-Text "95×56=" "83×93="

Replace-Text "42×27=" "42×20="
Replace-Text "61×50=" "33×30="
Replace-Text "76×51=" "80×22="
Replace-Text "51×78=" "41×21="
Replace-Text "14×12=" "23×59="
